# TDD CASE_1_ITER_like_LTS / transitory_input.xlsx
# Update the SIMULATION test-case-name formula to append "case1" to the
# TEXTJOIN'd string, and normalise the formatting of a handful of data-entry
# cells (E4, B5, C5, E5, E6, E7, E8, E10, E11) so they share the same
# centered / unlocked style as the rest of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TRANSIENT")

# --- E3: extend the TEXTJOIN formula with an extra "case1" literal ------
$ws.Range("E3").Formula = '=_xlfn.TEXTJOIN("_",TRUE,A6,E6,A8,E8,[1]GRID!$A$4,[1]GRID!$E$4,"case1")'

# --- Normalise formatting on the cells whose redundant style got merged -
$xlHAlignCenter = -4108

$cellsToRestyle = @("E3", "E4", "B5", "C5", "E5", "E6", "E7", "E8", "E10", "E11")
foreach ($addr in $cellsToRestyle) {
    $rng = $ws.Range($addr)
    $rng.HorizontalAlignment = $xlHAlignCenter
    $rng.Locked = $false
}

$wb.Save()
